$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.391.29"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "1.848.12"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6299"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.70%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07623"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2931"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.10%  "

$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07745"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("D12").Value = "1.869.05"
$ws.Range("E12").Value = "  -5.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.00001115"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.006"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6793"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("D17").Value = "2.123.48"
$ws.Range("E17").Value = "  -6.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.194"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").Value = "29.408.56"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.520"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1399"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.462"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.301"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05595"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.116"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.040"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.855"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.157"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7124"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.585"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "

$ws.Range("D38").Value = "1.243.08"
$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01806"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.775"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.409"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9019"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "

$ws.Range("E46").Value = "  +1.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.161"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4016"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.689"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.977"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1121"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.51%  "
